$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add the new "Multi_qubit_fidelity" worksheet after the last existing sheet
# (i.e. after "Two_qubit_fidelity"), matching the target sheet order:
#   Porter_Thomas, Two_qubit_fidelity, Multi_qubit_fidelity
#
# A plain Worksheets.Add() drops the sheetPr/sheetFormatPr settings that the
# rest of the workbook uses, so instead we clone the last existing sheet
# (which carries the right <outlinePr>/<pageSetUpPr>/row-height metadata),
# move it to the end, wipe its old data and rename it.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Multi_qubit_fidelity"
$ws.Cells.Clear()

# Tab-data for the new sheet: 50 rows, tab-separated, columns A-F.
# Rows 26-50 only populate columns A-D (E/F intentionally blank there).
$tsv = @"
1.0	1.0000000000000004	1.0	1.0	1.0000000000000002	1.0000000000000004
1.0000000000000002	0.9999999999999996	1.0000000000000002	1.0	1.0	0.9999999999999998
1.0	1.0	0.9999999999999999	0.9999999999999998	1.0000000000000004	1.0
0.9999999999999999	0.9999999999999991	1.0	1.0000000000000004	1.0	0.9999999999999996
0.9975284716977462	0.9638618187874645	1.0000000000000002	1.0	1.0000000000000002	0.9999999999999998
0.9886388454407095	0.9351702982086929	0.9980287268179067	0.9728927079794467	0.9999999999999998	0.9999999999999991
0.9898317349735312	0.8446391903463212	0.9968024855988438	0.9559309917503636	0.9994706760021338	0.9976980298759844
0.9752778774485117	0.6797298850561199	0.996178515245844	0.9069514829111364	0.9994810045450907	0.9935039869289917
0.9793585488446043	0.4997808913296969	0.9924240007316799	0.8835650196835112	0.9985340266494855	0.9787215325688144
0.9595926891965265	0.3309510964961775	0.9898427624885022	0.8102793873831907	0.9984598640292467	0.9626687784873917
0.9569784556804083	0.24057316209785204	0.9789524027752247	0.7010148922266324	0.9965240211542609	0.926573427869234
0.9233610115943023	0.15216240754849636	0.9667113270674282	0.5952018560085253	0.9951811776189005	0.9005494853333628
0.9444920724405748	0.08689283163200259	0.9847098442774307	0.4285015225520468	0.9920057082791969	0.8661874063983098
0.9282494870049126	0.055702880937660196	0.9641259888566899	0.3595366099450983	0.9877736275722708	0.817292399559428
0.9436563314576419	0.03675471545489852	0.9662342440704986	0.23980954858793668	0.986461127199691	0.7219821447283694
0.9437942349767552	0.01905986532157093	0.9559262716045202	0.19271425781866994	0.9856306450608191	0.660323788453426
0.932750600935691	0.012178680423484928	0.9602689606768	0.12780697154610726	0.9820735828226811	0.55432353372269
0.9358548658681284	0.006342587645911301	0.9675735621164506	0.09012089206928639	0.9735773302212589	0.5061099111436245
0.9378564886355464	0.004239419577894983	0.9639240011424405	0.06465502716798047	0.9810189881824312	0.41539892347103013
0.947577661695712	0.0023039257039833495	0.9451615182325646	0.04659471910179421	0.9795135843340975	0.3404802474106149
0.9380232610630282	0.0012348000468564593	0.9600414049450341	0.031596485932728326	0.9682638940481276	0.26482609037987626
0.9381841340426768	0.000584313636456046	0.9637391505047291	0.02075789254381546	0.9845981615245205	0.20661586783339625
0.9379540709207703	0.0003859948893576825	0.9646936240892926	0.014013221050685504	0.966818487963842	0.16838602122445884
0.9332787530591284	0.00022350476698166467	0.9706278109850908	0.009470060404444761	0.9712022836852451	0.12047195409154157
0.9456431931321851	0.000131351230917407	0.9743212240703274	0.006384432808249788	0.9705174357018557	0.08869614487003318
0.9441062250730744	7.474895266518346e-05	0.9566798774864592	0.004101749378490087		
0.9487582662460302	3.789046155974221e-05	0.9712621987827051	0.0027668462661666537		
0.9249656829146858	2.7412970444028004e-05	0.9513517469039859	0.0018655497037977533		
0.933282441431512	1.861121578886998e-05	0.9647551627081123	0.0011137849128783509		
0.9336898090504143	1.0119273068577818e-05	0.9542376917583547	0.0007428216115198193		
0.9426443473547501	1.115755948796743e-05	0.9630758465041722	0.00040489026530680916		
0.9444895264879868	4.1321242158652866e-06	0.9508294525941141	0.0002441542661199678		
0.9408895838427395	2.1420345069005633e-06	0.9594518100013029	0.00015726918759493668		
0.9291996948318093	5.078482219863495e-06	0.9482578711806148	8.548426837186698e-05		
0.9590866750507065	4.458532235532586e-06	0.9648647682165917	5.536343341494502e-05		
0.9366034227621188	2.1502647877915326e-06	0.9670206113983673	2.3550338224417508e-05		
0.9505919810170279	1.7812273055186353e-06	0.9583038545147248	1.4045074089349634e-05		
0.9099030499953744	7.844469463225438e-07	0.9586479369127683	8.167610419351914e-06		
0.9552456676086443	1.2040415531779482e-06	0.9579954807400195	7.2369687035348445e-06		
0.9337243671909122	1.0166942654508274e-07	0.9649651689795758	8.039775332109232e-06		
0.9517759566033444	4.9892300091114114e-06	0.9678765441729926	9.202078131289957e-06		
0.9398401428286413	3.404004836572968e-06	0.9631047556049804	3.885100070422563e-06		
0.946190174891669	1.3233895622277472e-06	0.9617277298127936	4.321388599904069e-06		
0.9347172999646729	2.582298725781132e-07	0.9673216096637237	2.8666808434095752e-06		
0.9402610767692726	1.8618509288570674e-07	0.9565116473132251	1.9491412147595098e-06		
0.9413203493723766	8.027484336539393e-08	0.9689006658976205	1.3697626422234069e-06		
0.9363302582951751	1.0912883757892506e-06	0.9544088866841088	4.872141896653974e-07		
0.9251761184383428	1.277790491093836e-06	0.956241726707407	1.5711523110528523e-06		
0.9480978453950172	8.450937189286067e-07	0.9625595216463121	2.895788350463419e-07		
0.9318701395874858	5.273613983788518e-07	0.9478193435603651	6.326773890466455e-07		
"@

$lines = $tsv -split "`r?`n"
$rowNum = 0
foreach ($line in $lines) {
    $rowNum++
    if ($rowNum -gt 50) { break }
    $fields = $line -split "`t"
    for ($colIdx = 0; $colIdx -lt $fields.Length; $colIdx++) {
        $text = $fields[$colIdx]
        if ($text -ne "") {
            $cell = $ws.Cells.Item($rowNum, $colIdx + 1)
            # Force text storage (the source data are numeric-looking strings,
            # not numbers) and then strip the temporary formatting so the
            # cell is left with the default (unstyled) look.
            $cell.NumberFormat = "@"
            $cell.Value = $text
            $cell.ClearFormats()
        }
    }
}

# Selection on the new sheet should rest at A1 (it becomes the active sheet,
# so Excel also marks its tab as selected automatically).
$ws.Range("A1").Select()
